$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "中芯国际"
$ws.Range("B2").Value = "中芯国际"
$ws.Range("C2").Value = "中芯国际"
$ws.Range("A3").Value = "赣锋锂业"
$ws.Range("B3").Value = "赣锋锂业"
$ws.Range("A4").Value = "海南华铁"
$ws.Range("B4").Value = "天齐锂业"
$ws.Range("C4").Value = "海南华铁"
$ws.Range("A5").Value = "江波龙"
$ws.Range("B5").Value = "海南华铁"
$ws.Range("C5").Value = "赣锋锂业"
$ws.Range("A6").Value = "山子高科"
$ws.Range("C6").Value = "三花智控"
$ws.Range("A7").Value = "深科技"
$ws.Range("B7").Value = "东方财富"
$ws.Range("C7").Value = "张江高科"
$ws.Range("A8").Value = "天齐锂业"
$ws.Range("B8").Value = "山子高科"
$ws.Range("A9").Value = "三花智控"
$ws.Range("C9").Value = "赛力斯"
$ws.Range("A10").Value = "张江高科"
$ws.Range("B10").Value = "华友钴业"
$ws.Range("C10").Value = "紫金矿业"
$ws.Range("A11").Value = "宁德时代"
$ws.Range("B11").Value = "宁德时代"
$ws.Range("C11").Value = "天齐锂业"
$ws.Range("A12").Value = "XD紫金矿"
$ws.Range("B12").Value = "贵州茅台"
$ws.Range("C12").Value = "天赐材料"
$ws.Range("A13").Value = "长电科技"
$ws.Range("B13").Value = "江波龙"
$ws.Range("C13").Value = "万向钱潮"
$ws.Range("A14").Value = "华友钴业"
$ws.Range("B14").Value = "长电科技"
$ws.Range("C14").Value = "长电科技"
$ws.Range("A15").Value = "天赐材料"
$ws.Range("B15").Value = "张江高科"
$ws.Range("C15").Value = "中电鑫龙"
$ws.Range("A16").Value = "德明利"
$ws.Range("B16").Value = "洛阳钼业"
$ws.Range("C16").Value = "大众公用"
$ws.Range("A17").Value = "东方财富"
$ws.Range("B17").Value = "XD紫金矿"
$ws.Range("C17").Value = "领益智造"
$ws.Range("A18").Value = "中信证券"
$ws.Range("B18").Value = "江西铜业"
$ws.Range("C18").Value = "宁德时代"
$ws.Range("A19").Value = "上海电气"
$ws.Range("B19").Value = "上海电气"
$ws.Range("C19").Value = "华友钴业"
$ws.Range("A20").Value = "平安银行"
$ws.Range("B20").Value = "天赐材料"
$ws.Range("C20").Value = "立讯精密"
$ws.Range("B21").Value = "华虹公司"
$ws.Range("C21").Value = "上海电气"
